# Refactor phrase type determination logic
# Update the "PhraseType" column (G) so that all rows previously marked
# as "phrase" are now marked as "word", matching the unified phrase
# type determination logic described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Determine the last used row of the sheet.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($current -eq "phrase") {
        $cell.Value = "word"
    }
}
